# Update "想去人数" (want-to-go count) figures and mark two sold-out
# events as "不可售" (not for sale) in their "最低票价" (lowest price)
# column, across the 展览 / 演出 / 全部类型 sheets.
# Output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet 1) ----
$ws1.Range("F6").Value = 1163
$ws1.Range("F7").Value = 1592
$ws1.Range("F8").Value = 177
$ws1.Range("F9").Value = 177
$ws1.Range("F10").Value = 24
$ws1.Range("F11").Value = 1543
$ws1.Range("F12").Value = 3150
$ws1.Range("F13").Value = 690
$ws1.Range("F14").Value = 1857
$ws1.Range("F15").Value = 1837
$ws1.Range("F16").Value = 897
$ws1.Range("F17").Value = 308
$ws1.Range("F19").Value = 1520
$ws1.Range("F23").Value = 1311
$ws1.Range("F24").Value = 428
$ws1.Range("F25").Value = 528
$ws1.Range("F26").Value = 209
$ws1.Range("F27").Value = 7787
$ws1.Range("F28").Value = 7206
$ws1.Range("F29").Value = 780
$ws1.Range("F31").Value = 1727
$ws1.Range("F32").Value = 102
$ws1.Range("F33").Value = 256

# ---- 演出 (sheet 2) ----
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F7").Value = 28
$ws2.Range("F8").Value = 3

# ---- 全部类型 (sheet 4) ----
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F8").Value = 1163
$ws4.Range("F9").Value = 1592
$ws4.Range("F10").Value = 177
$ws4.Range("F11").Value = 177
$ws4.Range("F13").Value = 24
$ws4.Range("F14").Value = 1543
$ws4.Range("F15").Value = 3150
$ws4.Range("F16").Value = 690
$ws4.Range("F17").Value = 1857
$ws4.Range("F18").Value = 1837
$ws4.Range("F19").Value = 897
$ws4.Range("F20").Value = 308
$ws4.Range("F22").Value = 1520
$ws4.Range("F28").Value = 1311
$ws4.Range("F29").Value = 428
$ws4.Range("F30").Value = 528
$ws4.Range("F31").Value = 209
$ws4.Range("F32").Value = 7787
$ws4.Range("F33").Value = 7206
$ws4.Range("F34").Value = 781
$ws4.Range("F36").Value = 1727
$ws4.Range("F38").Value = 28
$ws4.Range("F39").Value = 102
$ws4.Range("F40").Value = 256
$ws4.Range("F41").Value = 3
